$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "hello"
Write-Output ("A1=" + $ws.Range("A1").Value2)
$ws.Range("B1").Value2 = "world"
Write-Output ("B1=" + $ws.Range("B1").Value2)
